# Updated cryptos list values (Price / Volume(1h)) per upstream data refresh.
# Applies per-row D (Price) and E (Volume(1h)) text updates to sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.191.45"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.549.26"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "592.02"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "173.12"
$ws.Range("E6").Value = "  +4.84%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "2.548.89"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "0.349"
$ws.Range("E13").Value = "  -4.72%  "
$ws.Range("D14").Value = "27.05"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "3.017.45"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "67.049.32"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "2.552.45"
$ws.Range("E18").Value = "  -2.56%  "
$ws.Range("D19").Value = "7.99"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").Value = "11.37"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "355.84"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("E24").Value = "  +5.71%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "70.03"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").Value = "2.680.18"
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "537.62"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").Value = "8.18"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "158.51"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "18.73"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "18.45"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "5.17"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").Value = "2.53"
$ws.Range("E44").Value = "  +5.30%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "39.72"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").Value = "150.48"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "0.564"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("D49").Value = "0.0₆0281"
$ws.Range("E49").Value = "  -5.22%  "
$ws.Range("D50").Value = "3.71"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("E51").Value = "  +0.57%  "
